$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.028.47'
$ws.Range('E2').Value = '  -1.70%  '
$ws.Range('D3').Value = '1.553.56'
$ws.Range('E3').Value = '  -0.90%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.000'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '287.67'
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3981'
$ws.Range('E7').Value = '  +6.25%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3229'
$ws.Range('E8').Value = '  -1.43%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.61'
$ws.Range('E9').Value = '  -6.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07327'
$ws.Range('E10').Value = '  -1.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.099'
$ws.Range('E11').Value = '  -4.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.91'
$ws.Range('E13').Value = '  -7.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.681'
$ws.Range('E14').Value = '  -2.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.724'
$ws.Range('E15').Value = '  -1.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001129'
$ws.Range('E16').Value = '  +2.61%  '
$ws.Range('D17').Value = '1.555.82'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06615'
$ws.Range('E18').Value = '  -1.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '85.08'
$ws.Range('E19').Value = '  -0.92%  '
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.334'
$ws.Range('E21').Value = '  -0.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.87'
$ws.Range('E22').Value = '  -2.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.27'
$ws.Range('E23').Value = '  -3.75%  '
$ws.Range('D24').Value = '22.027.53'
$ws.Range('E24').Value = '  -1.74%  '
$ws.Range('E25').Value = '  +0.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.461'
$ws.Range('E26').Value = '  -3.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '148.02'
$ws.Range('E28').Value = '  -3.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.860'
$ws.Range('E29').Value = '  -1.18%  '
$ws.Range('D30').Value = '1.718.51'
$ws.Range('E30').Value = '  -1.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '120.84'
$ws.Range('E31').Value = '  -2.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.066'
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.722'
$ws.Range('E33').Value = '  -3.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08393'
$ws.Range('E34').Value = '  +1.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.277'
$ws.Range('E35').Value = '  -3.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.629'
$ws.Range('E36').Value = '  -16.13%  '
$ws.Range('E37').Value = '  -1.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02275'
$ws.Range('E38').Value = '  -4.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.127'
$ws.Range('E39').Value = '  -2.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2097'
$ws.Range('E40').Value = '  -4.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.212'
$ws.Range('E41').Value = '  -6.75%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.66'
$ws.Range('E43').Value = '  -4.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5871'
$ws.Range('E44').Value = '  -3.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.26'
$ws.Range('E45').Value = '  -3.85%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.712'
$ws.Range('E46').Value = '  -1.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5602'
$ws.Range('E47').Value = '  -5.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.912'
$ws.Range('E48').Value = '  -4.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '117.81'
$ws.Range('E49').Value = '  -4.80%  '
$ws.Range('E50').Value = '  -2.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06857'
$ws.Range('E51').Value = '  -4.10%  '
